$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.750068724155426
$ws.Range("B1").Value = 1.157423138618469
$ws.Range("C1").Value = 3.547337055206299
$ws.Range("D1").Value = 2.257123231887817
$ws.Range("E1").Value = 0.8098844885826111
